$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Antal" (I) column keeps storing its numeric-looking counts as text,
# matching the source workbook where I2:I10 are inline strings (not numbers).
$ws.Range("I2:I10").NumberFormat = "@"

# Step 1: capture the current ("before") values of the columns that move between rows
$r2_A = $ws.Cells.Item(2,1).Value2
$r2_B = $ws.Cells.Item(2,2).Value2
$r2_D = $ws.Cells.Item(2,4).Value2
$r2_E = $ws.Cells.Item(2,5).Value2
$r2_F = $ws.Cells.Item(2,6).Value2
$r2_G = $ws.Cells.Item(2,7).Value2
$r2_H = $ws.Cells.Item(2,8).Value2
$r2_I = $ws.Cells.Item(2,9).Value2
$r2_Q = $ws.Cells.Item(2,17).Value2
$r2_R = $ws.Cells.Item(2,18).Value2
$r2_AO = $ws.Cells.Item(2,41).Value2
$r3_A = $ws.Cells.Item(3,1).Value2
$r3_B = $ws.Cells.Item(3,2).Value2
$r3_D = $ws.Cells.Item(3,4).Value2
$r3_E = $ws.Cells.Item(3,5).Value2
$r3_F = $ws.Cells.Item(3,6).Value2
$r3_G = $ws.Cells.Item(3,7).Value2
$r3_H = $ws.Cells.Item(3,8).Value2
$r3_I = $ws.Cells.Item(3,9).Value2
$r3_Q = $ws.Cells.Item(3,17).Value2
$r3_R = $ws.Cells.Item(3,18).Value2
$r3_AO = $ws.Cells.Item(3,41).Value2
$r4_A = $ws.Cells.Item(4,1).Value2
$r4_B = $ws.Cells.Item(4,2).Value2
$r4_D = $ws.Cells.Item(4,4).Value2
$r4_E = $ws.Cells.Item(4,5).Value2
$r4_F = $ws.Cells.Item(4,6).Value2
$r4_G = $ws.Cells.Item(4,7).Value2
$r4_H = $ws.Cells.Item(4,8).Value2
$r4_I = $ws.Cells.Item(4,9).Value2
$r4_Q = $ws.Cells.Item(4,17).Value2
$r4_R = $ws.Cells.Item(4,18).Value2
$r4_AO = $ws.Cells.Item(4,41).Value2
$r5_A = $ws.Cells.Item(5,1).Value2
$r5_B = $ws.Cells.Item(5,2).Value2
$r5_D = $ws.Cells.Item(5,4).Value2
$r5_E = $ws.Cells.Item(5,5).Value2
$r5_F = $ws.Cells.Item(5,6).Value2
$r5_G = $ws.Cells.Item(5,7).Value2
$r5_H = $ws.Cells.Item(5,8).Value2
$r5_I = $ws.Cells.Item(5,9).Value2
$r5_Q = $ws.Cells.Item(5,17).Value2
$r5_R = $ws.Cells.Item(5,18).Value2
$r5_AO = $ws.Cells.Item(5,41).Value2
$r6_A = $ws.Cells.Item(6,1).Value2
$r6_B = $ws.Cells.Item(6,2).Value2
$r6_D = $ws.Cells.Item(6,4).Value2
$r6_E = $ws.Cells.Item(6,5).Value2
$r6_F = $ws.Cells.Item(6,6).Value2
$r6_G = $ws.Cells.Item(6,7).Value2
$r6_H = $ws.Cells.Item(6,8).Value2
$r6_I = $ws.Cells.Item(6,9).Value2
$r6_Q = $ws.Cells.Item(6,17).Value2
$r6_R = $ws.Cells.Item(6,18).Value2
$r6_AO = $ws.Cells.Item(6,41).Value2
$r7_A = $ws.Cells.Item(7,1).Value2
$r7_B = $ws.Cells.Item(7,2).Value2
$r7_D = $ws.Cells.Item(7,4).Value2
$r7_E = $ws.Cells.Item(7,5).Value2
$r7_F = $ws.Cells.Item(7,6).Value2
$r7_G = $ws.Cells.Item(7,7).Value2
$r7_H = $ws.Cells.Item(7,8).Value2
$r7_I = $ws.Cells.Item(7,9).Value2
$r7_Q = $ws.Cells.Item(7,17).Value2
$r7_R = $ws.Cells.Item(7,18).Value2
$r7_AO = $ws.Cells.Item(7,41).Value2
$r8_A = $ws.Cells.Item(8,1).Value2
$r8_B = $ws.Cells.Item(8,2).Value2
$r8_D = $ws.Cells.Item(8,4).Value2
$r8_E = $ws.Cells.Item(8,5).Value2
$r8_F = $ws.Cells.Item(8,6).Value2
$r8_G = $ws.Cells.Item(8,7).Value2
$r8_H = $ws.Cells.Item(8,8).Value2
$r8_I = $ws.Cells.Item(8,9).Value2
$r8_Q = $ws.Cells.Item(8,17).Value2
$r8_R = $ws.Cells.Item(8,18).Value2
$r8_AO = $ws.Cells.Item(8,41).Value2
$r9_A = $ws.Cells.Item(9,1).Value2
$r9_B = $ws.Cells.Item(9,2).Value2
$r9_D = $ws.Cells.Item(9,4).Value2
$r9_E = $ws.Cells.Item(9,5).Value2
$r9_F = $ws.Cells.Item(9,6).Value2
$r9_G = $ws.Cells.Item(9,7).Value2
$r9_H = $ws.Cells.Item(9,8).Value2
$r9_I = $ws.Cells.Item(9,9).Value2
$r9_Q = $ws.Cells.Item(9,17).Value2
$r9_R = $ws.Cells.Item(9,18).Value2
$r9_AO = $ws.Cells.Item(9,41).Value2
$r10_A = $ws.Cells.Item(10,1).Value2
$r10_B = $ws.Cells.Item(10,2).Value2
$r10_D = $ws.Cells.Item(10,4).Value2
$r10_E = $ws.Cells.Item(10,5).Value2
$r10_F = $ws.Cells.Item(10,6).Value2
$r10_G = $ws.Cells.Item(10,7).Value2
$r10_H = $ws.Cells.Item(10,8).Value2
$r10_I = $ws.Cells.Item(10,9).Value2
$r10_Q = $ws.Cells.Item(10,17).Value2
$r10_R = $ws.Cells.Item(10,18).Value2
$r10_AO = $ws.Cells.Item(10,41).Value2

# Step 2: write the captured values into their new rows (rows 2-10 are permuted)
# row 2 gets the data that used to be in row 8
$ws.Cells.Item(2,1).Value = $r8_A
$ws.Cells.Item(2,2).Value = $r8_B
$ws.Cells.Item(2,4).Value = $r8_D
$ws.Cells.Item(2,5).Value = $r8_E
$ws.Cells.Item(2,6).Value = $r8_F
$ws.Cells.Item(2,7).Value = $r8_G
$ws.Cells.Item(2,8).Value = $r8_H
$ws.Cells.Item(2,9).Value = $r8_I
$ws.Cells.Item(2,17).Value = $r8_Q
$ws.Cells.Item(2,18).Value = $r8_R
$ws.Cells.Item(2,41).Value = $r8_AO

# row 3 gets the data that used to be in row 6
$ws.Cells.Item(3,1).Value = $r6_A
$ws.Cells.Item(3,2).Value = $r6_B
$ws.Cells.Item(3,4).Value = $r6_D
$ws.Cells.Item(3,5).Value = $r6_E
$ws.Cells.Item(3,6).Value = $r6_F
$ws.Cells.Item(3,7).Value = $r6_G
$ws.Cells.Item(3,8).Value = $r6_H
$ws.Cells.Item(3,9).Value = $r6_I
$ws.Cells.Item(3,17).Value = $r6_Q
$ws.Cells.Item(3,18).Value = $r6_R
$ws.Cells.Item(3,41).Value = $r6_AO

# row 4 gets the data that used to be in row 3
$ws.Cells.Item(4,1).Value = $r3_A
$ws.Cells.Item(4,2).Value = $r3_B
$ws.Cells.Item(4,4).Value = $r3_D
$ws.Cells.Item(4,5).Value = $r3_E
$ws.Cells.Item(4,6).Value = $r3_F
$ws.Cells.Item(4,7).Value = $r3_G
$ws.Cells.Item(4,8).Value = $r3_H
$ws.Cells.Item(4,9).Value = $r3_I
$ws.Cells.Item(4,17).Value = $r3_Q
$ws.Cells.Item(4,18).Value = $r3_R
$ws.Cells.Item(4,41).Value = $r3_AO

# row 5 gets the data that used to be in row 2
$ws.Cells.Item(5,1).Value = $r2_A
$ws.Cells.Item(5,2).Value = $r2_B
$ws.Cells.Item(5,4).Value = $r2_D
$ws.Cells.Item(5,5).Value = $r2_E
$ws.Cells.Item(5,6).Value = $r2_F
$ws.Cells.Item(5,7).Value = $r2_G
$ws.Cells.Item(5,8).Value = $r2_H
$ws.Cells.Item(5,9).Value = $r2_I
$ws.Cells.Item(5,17).Value = $r2_Q
$ws.Cells.Item(5,18).Value = $r2_R
$ws.Cells.Item(5,41).Value = $r2_AO

# row 6 gets the data that used to be in row 9
$ws.Cells.Item(6,1).Value = $r9_A
$ws.Cells.Item(6,2).Value = $r9_B
$ws.Cells.Item(6,4).Value = $r9_D
$ws.Cells.Item(6,5).Value = $r9_E
$ws.Cells.Item(6,6).Value = $r9_F
$ws.Cells.Item(6,7).Value = $r9_G
$ws.Cells.Item(6,8).Value = $r9_H
$ws.Cells.Item(6,9).Value = $r9_I
$ws.Cells.Item(6,17).Value = $r9_Q
$ws.Cells.Item(6,18).Value = $r9_R
$ws.Cells.Item(6,41).Value = $r9_AO

# row 7 gets the data that used to be in row 10
$ws.Cells.Item(7,1).Value = $r10_A
$ws.Cells.Item(7,2).Value = $r10_B
$ws.Cells.Item(7,4).Value = $r10_D
$ws.Cells.Item(7,5).Value = $r10_E
$ws.Cells.Item(7,6).Value = $r10_F
$ws.Cells.Item(7,7).Value = $r10_G
$ws.Cells.Item(7,8).Value = $r10_H
$ws.Cells.Item(7,9).Value = $r10_I
$ws.Cells.Item(7,17).Value = $r10_Q
$ws.Cells.Item(7,18).Value = $r10_R
$ws.Cells.Item(7,41).Value = $r10_AO

# row 8 gets the data that used to be in row 7
$ws.Cells.Item(8,1).Value = $r7_A
$ws.Cells.Item(8,2).Value = $r7_B
$ws.Cells.Item(8,4).Value = $r7_D
$ws.Cells.Item(8,5).Value = $r7_E
$ws.Cells.Item(8,6).Value = $r7_F
$ws.Cells.Item(8,7).Value = $r7_G
$ws.Cells.Item(8,8).Value = $r7_H
$ws.Cells.Item(8,9).Value = $r7_I
$ws.Cells.Item(8,17).Value = $r7_Q
$ws.Cells.Item(8,18).Value = $r7_R
$ws.Cells.Item(8,41).Value = $r7_AO

# row 9 gets the data that used to be in row 5
$ws.Cells.Item(9,1).Value = $r5_A
$ws.Cells.Item(9,2).Value = $r5_B
$ws.Cells.Item(9,4).Value = $r5_D
$ws.Cells.Item(9,5).Value = $r5_E
$ws.Cells.Item(9,6).Value = $r5_F
$ws.Cells.Item(9,7).Value = $r5_G
$ws.Cells.Item(9,8).Value = $r5_H
$ws.Cells.Item(9,9).Value = $r5_I
$ws.Cells.Item(9,17).Value = $r5_Q
$ws.Cells.Item(9,18).Value = $r5_R
$ws.Cells.Item(9,41).Value = $r5_AO

# row 10 gets the data that used to be in row 4
$ws.Cells.Item(10,1).Value = $r4_A
$ws.Cells.Item(10,2).Value = $r4_B
$ws.Cells.Item(10,4).Value = $r4_D
$ws.Cells.Item(10,5).Value = $r4_E
$ws.Cells.Item(10,6).Value = $r4_F
$ws.Cells.Item(10,7).Value = $r4_G
$ws.Cells.Item(10,8).Value = $r4_H
$ws.Cells.Item(10,9).Value = $r4_I
$ws.Cells.Item(10,17).Value = $r4_Q
$ws.Cells.Item(10,18).Value = $r4_R
$ws.Cells.Item(10,41).Value = $r4_AO
